$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RemoveTest")

# New data for column A (rows 2..51)
$colA = @(10845500,556400,412400,333900,1444500,183800,1464900,680000,3319300,332600,312400,322500,336200,3675700,2841700,2723200,1738300,1738800,177500,181500,187500,1363700,162500,152200,194100,183800,153000,151800,157200,154500,154700,156900,155500,154600,157300,156900,156500,156900,152700,156300,153800,155300,1344400,161700,156800,149000,1316700,157000,1356000,178000)

# New data for column B (rows 2..51)
$colB = @(107000,107000,107000,107000,107000,107000,107000,106971,107000,106971,107000,107000,107000,107000,107000,107000,107000,106979,106973,107000,107000,107000,107000,107000,106983,107000,107000,107000,106983,107000,107000,107000,107000,107000,107000,107000,106970,107000,106965,107000,107000,106971,106964,106969,107000,107000,107000,107000,107000,107000)

for ($i = 0; $i -lt $colA.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

# Update the selection to F4 to reflect the last active cell in the sheet
$ws.Range("F4").Select() | Out-Null

$wb.Save()
